# Revert "updata project members again"
# This reverts the addition of the "Longbo Qiao" row (row 11) and restores
# Joon Lee's (row 5) earlier contact e-mail address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Restore Joon Lee's (row 5) preferred contact e-mail back to the
#    earlier address.
$ws.Range("D5").Value = "ljshope@stanford.kr"

# 2. Remove the hyperlink that was attached to D11 (Longbo Qiao's e-mail)
#    before the row itself is deleted, so no dangling hyperlink reference
#    is left behind.
$linksToRemove = @()
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$11') {
        $linksToRemove += $hl
    }
}
foreach ($hl in $linksToRemove) {
    $hl.Delete()
}

# 3. Delete row 11 (the "Longbo Qiao" entry) entirely; this also shrinks
#    the sheet dimension back down to A1:F10 automatically.
$ws.Rows("11:11").Delete() | Out-Null

# 4. Restore the previously selected/active cell to E9 (it pointed at
#    E11 while the now-removed row existed).
$ws.Range("E9").Select() | Out-Null
